$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows just above the "Total hours Spent" block (old row 36),
# shifting it down to rows 39:41 and opening up rows 33:35 for new content.
$ws.Rows("33:35").Insert()

# The inserted rows picked up row 32's formatting; re-apply the per-column
# style pattern (A=plain, B=date, C=plain, D=plain) used by the rest of the
# table (e.g. row 9) across the 4 affected rows (33:36).
$ws.Range("A9:D9").Copy()
$ws.Range("A33:D36").PasteSpecial(-4122)

# New data row logged for the day.
$ws.Range("A33").Value = 28
$ws.Range("B33").Value = 45326
$ws.Range("C33").Value = 8
$ws.Range("D33").Value = "Improved security a lot"

# Update the view state recorded for the sheet (scrolled down, new active cell).
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$null = $ws.Range("B39").Select()
